$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.91016218945221
$ws.Range("C2").Value = 9.087496150431338
$ws.Range("D2").Value = 6.000319336270256
$ws.Range("E2").Value = 12.39933772938301
$ws.Range("F2").Value = 48.35214782205405
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.51387118435646
$ws.Range("K2").Value = 16.44594121733526
$ws.Range("N2").Value = 23.38362989782224
$ws.Range("B3").Value = 16.72091115277286
$ws.Range("C3").Value = 8.955958201609997
$ws.Range("D3").Value = 5.991839123878377
$ws.Range("E3").Value = 12.36723847085993
$ws.Range("F3").Value = 48.23118647083447
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.52161106110634
$ws.Range("K3").Value = 16.32495903335073
$ws.Range("N3").Value = 23.42070127150854
$ws.Range("B4").Value = 16.60853212556576
$ws.Range("C4").Value = 8.877108239591166
$ws.Range("D4").Value = 5.987983214063489
$ws.Range("E4").Value = 12.35036435372198
$ws.Range("F4").Value = 48.16678383780981
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.52809246315978
$ws.Range("K4").Value = 16.25463664351187
$ws.Range("N4").Value = 23.44533657101206
$ws.Range("B5").Value = 16.56375058004287
$ws.Range("C5").Value = 8.845500106932564
$ws.Range("D5").Value = 5.986753500837634
$ws.Range("E5").Value = 12.34420545500079
$ws.Range("F5").Value = 48.14303233604223
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.53116836135043
$ws.Range("K5").Value = 16.2270023697481
$ws.Range("N5").Value = 23.45584649757815
$ws.Range("B6").Value = 16.55637736940034
$ws.Range("C6").Value = 8.840284457531043
$ws.Range("D6").Value = 5.986569989859892
$ws.Range("E6").Value = 12.34322622546907
$ws.Range("F6").Value = 48.13923928813477
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.53170535948528
$ws.Range("K6").Value = 16.22247622628875
$ws.Range("N6").Value = 23.45762009422258
$ws.Range("B7").Value = 16.60792401301803
$ws.Range("C7").Value = 8.876679784977831
$ws.Range("D7").Value = 5.987965244282878
$ws.Range("E7").Value = 12.35027838226748
$ws.Range("F7").Value = 48.16645340883862
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.52813218606373
$ws.Range("K7").Value = 16.25425978358778
$ws.Range("N7").Value = 23.44547640549077
$ws.Range("B8").Value = 16.84414751895503
$ws.Range("C8").Value = 9.04177014692424
$ws.Range("D8").Value = 5.997116173144815
$ws.Range("E8").Value = 12.38768446487485
$ws.Range("F8").Value = 48.30839991894342
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.51618098100041
$ws.Range("K8").Value = 16.40342095570009
$ws.Range("N8").Value = 23.39602316469463
$ws.Range("B9").Value = 17.33528982096486
$ws.Range("C9").Value = 9.378803075830765
$ws.Range("D9").Value = 6.02568649149546
$ws.Range("E9").Value = 12.4832995369509
$ws.Range("F9").Value = 48.66444337776679
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.50646632548449
$ws.Range("K9").Value = 16.72609915938941
$ws.Range("N9").Value = 23.31391977522005
$ws.Range("B10").Value = 17.70980688872157
$ws.Range("C10").Value = 9.631949517185545
$ws.Range("D10").Value = 6.05301268193661
$ws.Range("E10").Value = 12.56676334762357
$ws.Range("F10").Value = 48.97242593026533
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.50769215369598
$ws.Range("K10").Value = 16.97977114404474
$ws.Range("N10").Value = 23.26268080188114
$ws.Range("B11").Value = 17.8824270502616
$ws.Range("C11").Value = 9.747769528388208
$ws.Range("D11").Value = 6.066784283537935
$ws.Range("E11").Value = 12.60751263761445
$ws.Range("F11").Value = 49.12235707454228
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.51006300540334
$ws.Range("K11").Value = 17.09836812998663
$ws.Range("N11").Value = 23.24134640986414
$ws.Range("B12").Value = 17.94805639747736
$ws.Range("C12").Value = 9.791678404488028
$ws.Range("D12").Value = 6.072188636688208
$ws.Range("E12").Value = 12.62333459381571
$ws.Range("F12").Value = 49.18052026700839
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.51122102945535
$ws.Range("K12").Value = 17.14370204458336
$ws.Range("N12").Value = 23.23355188286364
$ws.Range("B13").Value = 17.93391127624953
$ws.Range("C13").Value = 9.782220305544008
$ws.Range("D13").Value = 6.071016348840091
$ws.Range("E13").Value = 12.61990980625811
$ws.Range("F13").Value = 49.16793250549743
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.51096006372221
$ws.Range("K13").Value = 17.13392031247565
$ws.Range("N13").Value = 23.2352179206732
$ws.Range("B14").Value = 17.8878215206382
$ws.Range("C14").Value = 9.75138118201475
$ws.Range("D14").Value = 6.067225128771594
$ws.Range("E14").Value = 12.6088065344133
$ws.Range("F14").Value = 49.12711454263047
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.51015306462271
$ws.Range("K14").Value = 17.10208949550869
$ws.Range("N14").Value = 23.24069944833721
$ws.Range("B15").Value = 17.85962248818994
$ws.Range("C15").Value = 9.73249657254285
$ws.Range("D15").Value = 6.064927448913033
$ws.Range("E15").Value = 12.60205611748121
$ws.Range("F15").Value = 49.10229225651491
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.50969262728673
$ws.Range("K15").Value = 17.08264630173422
$ws.Range("N15").Value = 23.2440940848445
$ws.Range("B16").Value = 17.69856575014217
$ws.Range("C16").Value = 9.624389872870303
$ws.Range("D16").Value = 6.052139365310369
$ws.Range("E16").Value = 12.56415550821243
$ws.Range("F16").Value = 48.96282319393158
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.50757365942189
$ws.Range("K16").Value = 16.97208185633697
$ws.Range("N16").Value = 23.2641148215038
$ws.Range("B17").Value = 17.6002949872343
$ws.Range("C17").Value = 9.558207285992657
$ws.Range("D17").Value = 6.044635313179211
$ws.Range("E17").Value = 12.54161094466183
$ws.Range("F17").Value = 48.87976305780801
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.50673783262481
$ws.Range("K17").Value = 16.90504817470851
$ws.Range("N17").Value = 23.27690289916465
$ws.Range("B18").Value = 17.5439877532189
$ws.Range("C18").Value = 9.520205789552222
$ws.Range("D18").Value = 6.040445603762368
$ws.Range("E18").Value = 12.52890625640407
$ws.Range("F18").Value = 48.83291583705706
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.50642779425615
$ws.Range("K18").Value = 16.86679563584602
$ws.Range("N18").Value = 23.28444408460119
$ws.Range("B19").Value = 17.52496199388888
$ws.Range("C19").Value = 9.507351681802517
$ws.Range("D19").Value = 6.039048854786421
$ws.Range("E19").Value = 12.52464998514594
$ws.Range("F19").Value = 48.81721409720955
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.50635215288358
$ws.Range("K19").Value = 16.8538972001176
$ws.Range("N19").Value = 23.28702930280588
$ws.Range("B20").Value = 17.6107342222305
$ws.Range("C20").Value = 9.565246139335416
$ws.Range("D20").Value = 6.045421072335436
$ws.Range("E20").Value = 12.54398376199496
$ws.Range("F20").Value = 48.88850921390541
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.5068091435344
$ws.Range("K20").Value = 16.91215287885784
$ws.Range("N20").Value = 23.27552235178918
$ws.Range("B21").Value = 17.90135258277972
$ws.Range("C21").Value = 9.760438365735153
$ws.Range("D21").Value = 6.068333591908467
$ws.Range("E21").Value = 12.61205729413118
$ws.Range("F21").Value = 49.13906632682672
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.51038304182175
$ws.Range("K21").Value = 17.11142776935702
$ws.Range("N21").Value = 23.23908167039449
$ws.Range("B22").Value = 18.09278644859242
$ws.Range("C22").Value = 9.888281302133564
$ws.Range("D22").Value = 6.084410137210919
$ws.Range("E22").Value = 12.65882229310945
$ws.Range("F22").Value = 49.31089398622626
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.51423530597921
$ws.Range("K22").Value = 17.24411840808662
$ws.Range("N22").Value = 23.21692315420672
$ws.Range("B23").Value = 17.99049779171616
$ws.Range("C23").Value = 9.820038615469757
$ws.Range("D23").Value = 6.07573016159196
$ws.Range("E23").Value = 12.63365784886306
$ws.Range("F23").Value = 49.21845661165704
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.51204071974767
$ws.Range("K23").Value = 17.17308665647814
$ws.Range("N23").Value = 23.22859776284896
$ws.Range("B24").Value = 17.60601404774529
$ws.Range("C24").Value = 9.562063721566991
$ws.Range("D24").Value = 6.045065442663589
$ws.Range("E24").Value = 12.54291021172872
$ws.Range("F24").Value = 48.88455225527004
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.50677637278747
$ws.Range("K24").Value = 16.9089399484372
$ws.Range("N24").Value = 23.27614590806649
$ws.Range("B25").Value = 17.19977572108241
$ws.Range("C25").Value = 9.286469628901608
$ws.Range("D25").Value = 6.016833572396426
$ws.Range("E25").Value = 12.45508380190046
$ws.Range("F25").Value = 48.55989729316201
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.50762510545842
$ws.Range("K25").Value = 16.63575545616393
$ws.Range("N25").Value = 23.33453704752518
